$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet logs weekly price snapshots, one row per observation, with the
# most-recent week always inserted at row 2 (row 1 is the header) and older
# rows pushed down. This commit adds this week's new snapshot at row 2 and
# shifts every existing row down by one (old row 37 -> new row 38).

$lastRow = $ws.UsedRange.Rows.Count()
$newLastRow = $lastRow + 1

# Shift existing data rows (2..37) down by one row (38..3), working from the
# bottom up so we never overwrite a source row before it has been read.
for ($r = $lastRow; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":R" + $dstRow)
    $dst.Value2 = $src.Value2()
}

# The last row used to not exist, so it needs the date-column number format
# copied explicitly (plain Value2 writes don't touch formatting).
$ws.Range("D" + $newLastRow).NumberFormat = $ws.Range("D" + $lastRow).NumberFormat()

# Write this week's new snapshot into row 2.
$ws.Range("A2").Value2 = 11
$ws.Range("B2").Value2 = "Vega Monumental Concepción"
$ws.Range("C2").Value2 = "Bíobío"
$ws.Range("D2").Value2 = 44699
$ws.Range("E2").Value2 = 8
$ws.Range("F2").Value2 = 100112013
$ws.Range("G2").Value2 = "Alcachofa"
$ws.Range("H2").Value2 = "Española"
$ws.Range("I2").Value2 = "Primera"
$ws.Range("J2").Value2 = 100
$ws.Range("K2").Value2 = 19000
$ws.Range("L2").Value2 = 20000
$ws.Range("M2").Value2 = 19500
$ws.Range("N2").Value2 = "$/caja 30 unidades"
$ws.Range("O2").Value2 = "Provincia de Limarí"
$ws.Range("P2").Value2 = 650
$ws.Range("Q2").Value2 = 30
$ws.Range("R2").Value2 = "Hortaliza"

Write-Output "done"
